# Change aspect (column O, "ASP") to abs(180 - ASP) for all data rows.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.Cells.Item($ws.Rows.Count, 15).End(-4162).Row  # xlUp = -4162, column O = 15

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 15)
    $old = $cell.Value()
    if ($old -ne $null) {
        $new = [Math]::Abs(180 - $old)
        # Round-trip through ToString() so the stored double matches Excel's
        # native 15-significant-digit precision instead of raw IEEE-754 math.
        $cell.Value() = $new.ToString()
    }
}
